# Updates cryptos.xlsx price/volume figures (and a few reordered rows) to
# match the latest GitHub Actions scrape. All Price (D) / Volume(1h) (E)
# columns are plain text in this sheet, so numeric-looking price values
# are forced to text (NumberFormat "@") before assignment to avoid Excel
# silently converting them to numbers (which would also corrupt values
# like "11.00" / "1.00" by dropping trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.865.83'
$ws.Range("E2").Value = '  -1.24%  '
$ws.Range("D3").Value = '3.505.21'
$ws.Range("E3").Value = '  -1.56%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.68'
$ws.Range("E5").Value = '  +3.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '192.02'
$ws.Range("E6").Value = '  +1.86%  '
$ws.Range("E7").Value = '  +0.78%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("E10").Value = '  +3.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.50'
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.63'
$ws.Range("E13").Value = '  +2.29%  '
$ws.Range("D14").Value = '4.062.57'
$ws.Range("E14").Value = '  -1.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '620.96'
$ws.Range("E15").Value = '  +10.46%  '
$ws.Range("D16").Value = '69.928.41'
$ws.Range("E16").Value = '  -1.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.72'
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").Value = '3.498.55'
$ws.Range("E19").Value = '  -1.24%  '
$ws.Range("E20").Value = '  -0.06%  '
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.68'
$ws.Range("E22").Value = '  -1.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '105.85'
$ws.Range("E23").Value = '  +12.74%  '
$ws.Range("E24").Value = '  +0.75%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.01'
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.03'
$ws.Range("E26").Value = '  +4.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.00'
$ws.Range("E27").Value = '  -0.72%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.83'
$ws.Range("E28").Value = '  +5.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.17'
$ws.Range("E29").Value = '  +5.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.12'
$ws.Range("E30").Value = '  +1.09%  '
$ws.Range("B31").Value = 'dogwifhat'
$ws.Range("C31").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.17'
$ws.Range("E31").Value = '  +5.25%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.61'
$ws.Range("E32").Value = '  +3.52%  '
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '64.19'
$ws.Range("E34").Value = '  +1.70%  '
$ws.Range("D35").Value = '3.721.67'
$ws.Range("E35").Value = '  +2.10%  '
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '528.53'
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.10'
$ws.Range("E37").Value = '  -4.30%  '
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("E39").Value = '  +1.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.391'
$ws.Range("E40").Value = '  -3.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.80'
$ws.Range("E41").Value = '  -3.12%  '
$ws.Range("E42").Value = '  +0.82%  '
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0463'
$ws.Range("E44").Value = '  +1.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.86'
$ws.Range("E45").Value = '  -1.78%  '
$ws.Range("E46").Value = '  +2.56%  '
$ws.Range("E47").Value = '  -3.96%  '
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.75'
$ws.Range("E48").Value = '  -4.58%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("E49").Value = '  +0.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.03'
$ws.Range("E50").Value = '  -0.99%  '
$ws.Range("B51").Value = 'OceanProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.35'
$ws.Range("E51").Value = '  -6.41%  '
